$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2023-09-23 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-24 Sunday", 2) | Out-Null

$t = $d.Tables(1)

$t.Cell(1,1).Range.Text = "68÷6=11, 2"  # was: 84÷6=14, 0
$t.Cell(1,2).Range.Text = "95÷4=23, 3"  # was: 83÷9=9, 2
$t.Cell(1,3).Range.Text = "84÷7=12, 0"  # was: 61÷2=30, 1
$t.Cell(1,4).Range.Text = "21÷3=7, 0"  # was: 14÷4=3, 2
$t.Cell(1,5).Range.Text = "17÷5=3, 2"  # was: 25÷2=12, 1
$t.Cell(5,1).Range.Text = "84÷6=14, 0"  # was: 39÷9=4, 3
$t.Cell(5,2).Range.Text = "91÷2=45, 1"  # was: 68÷7=9, 5
$t.Cell(5,3).Range.Text = "14÷3=4, 2"  # was: 52÷9=5, 7
$t.Cell(5,4).Range.Text = "89÷2=44, 1"  # was: 31÷5=6, 1
$t.Cell(5,5).Range.Text = "44÷9=4, 8"  # was: 99÷4=24, 3
$t.Cell(9,1).Range.Text = "71÷4=17, 3"  # was: 11÷8=1, 3
$t.Cell(9,2).Range.Text = "10÷3=3, 1"  # was: 14÷4=3, 2
$t.Cell(9,3).Range.Text = "15÷6=2, 3"  # was: 57÷6=9, 3
$t.Cell(9,4).Range.Text = "20÷8=2, 4"  # was: 56÷8=7, 0
$t.Cell(9,5).Range.Text = "69÷2=34, 1"  # was: 56÷8=7, 0
$t.Cell(13,1).Range.Text = "40÷9=4, 4"  # was: 59÷2=29, 1
$t.Cell(13,2).Range.Text = "40÷8=5, 0"  # was: 25÷5=5, 0
$t.Cell(13,3).Range.Text = "19÷2=9, 1"  # was: 27÷7=3, 6
$t.Cell(13,4).Range.Text = "78÷9=8, 6"  # was: 41÷8=5, 1
$t.Cell(13,5).Range.Text = "82÷8=10, 2"  # was: 15÷5=3, 0
$t.Cell(17,1).Range.Text = "54÷9=6, 0"  # was: 81÷7=11, 4
$t.Cell(17,2).Range.Text = "10÷6=1, 4"  # was: 69÷3=23, 0
$t.Cell(17,3).Range.Text = "24÷5=4, 4"  # was: 53÷9=5, 8
$t.Cell(17,4).Range.Text = "17÷8=2, 1"  # was: 71÷7=10, 1
$t.Cell(17,5).Range.Text = "87÷9=9, 6"  # was: 58÷8=7, 2
